$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Holly added S.GISH to harvester in bioSamples -- update the harvester column (B)
# for all data rows to "S.GISH" (previously "Retrofitted_159").
$ws.Range("B2:B5").Value = "S.GISH"

# Reflect that column B (harvester) was the focus of this edit.
$ws.Columns("B:B").Select() | Out-Null
